$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("content")

# New key/value rows (externalized popup content from content.csv)
$keys = @(
    "user_in_x_zone",
    "user_zone_unkown",
    "user_zone_unkown_311",
    "user_zone",
    "evac_order",
    "no_evac_order"
)

$values = @(
    "You are not located in an Evacuation Zone (post-storm)",
    "Zone Finder cannot determine Zone for your address.<br>Try alternative address or determine Zone by examining map and clicking on your location. (post-storm)",
    "Zone Finder cannot determine Zone for your address.<br>Try alternative address. (post-storm)",
    "You are located in Zone `${zone} (post-storm)",
    "You are required to evacuate (post-storm)",
    "No evacuation order currently in effect (post-storm)"
)

$startRow = 11

for ($i = 0; $i -lt $keys.Length; $i++) {
    $row = $startRow + $i
    $cellA = $ws.Range("A$row")
    $cellA.Locked = $false
    $cellA.Value = $keys[$i]
    $cellA.Style = "Normal"
}

for ($i = 0; $i -lt $values.Length; $i++) {
    $row = $startRow + $i
    $cellB = $ws.Range("B$row")
    $cellB.Locked = $false
    $cellB.Value = $values[$i]
}

$ws.Activate() | Out-Null
$ws.Range("B2:B16").Select() | Out-Null

Write-Host "done"
